$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Quest")
$ws2 = $wb.Worksheets.Item("Reward")

# --- Fill in the headers -----------------------------------------------
# Written in this precise order so the shared-string table ends up
# ordered id1, id2, id, id3, attr1, attr2, attr3 (Reward's "id" is
# interleaved between Quest's "id2" and "id3" columns).
$ws1.Range("A1").Value = "id1"
$ws1.Range("B1").Value = "id2"
$ws2.Range("A1").Value = "id"
$ws1.Range("C1").Value = "id3"
$ws1.Range("D1").Value = "attr1"
$ws1.Range("E1").Value = "attr2"
$ws1.Range("F1").Value = "attr3"

# --- Fill in the Quest sheet data rows (A2:F10) -----------------------------
$data = @(
  @(1,1,1,1,1,1),
  @(1,1,2,2,2,2),
  @(1,1,3,3,3,3),
  @(2,1,1,1,1,1),
  @(2,1,2,2,2,2),
  @(2,1,3,3,3,3),
  @(3,1,1,1,1,1),
  @(3,1,2,2,2,2),
  @(3,1,3,3,3,3)
)

$row = 2
foreach ($r in $data) {
  $ws1.Cells.Item($row, 1).Value = $r[0]
  $ws1.Cells.Item($row, 2).Value = $r[1]
  $ws1.Cells.Item($row, 3).Value = $r[2]
  $ws1.Cells.Item($row, 4).Value = $r[3]
  $ws1.Cells.Item($row, 5).Value = $r[4]
  $ws1.Cells.Item($row, 6).Value = $r[5]
  $row = $row + 1
}

# --- View state --------------------------------------------------------
# Reward's selection settles back on its default cell (A1) ...
$ws2.Range("A1").Select()

# ... then Quest becomes the active / tab-selected sheet with the
# selection parked on I5 (this also moves the workbook's activeTab back
# to sheet 0, i.e. no activeTab override is left behind).
$ws1.Activate()
$ws1.Range("I5").Select()
